# "style: Restyled start screen" - leaderboard data grew (more rows recorded
# for Sheet1/Sheet4/Sheet5), and two new player-name strings ("p", "q") were
# introduced via the new Sheet1 rows.

$wb = $excel.ActiveWorkbook

# ---- Sheet1 ("Sheet1", rId1 -> sheet1.xml): append rows 22-33 ----
$ws1 = $wb.Worksheets.Item("Sheet1")
$sheet1Rows = @(
    @("Jack", 521.0),
    @("Jack", 0.0),
    @("Jack", 1563.0),
    @("p",    1042.0),
    @("a",    4168.0),
    @("q",    1042.0),
    @("a",    1042.0),
    @("a",    1042.0),
    @("j",    1042.0),
    @("l",    1042.0),
    @("a",    1042.0),
    @("a",    0.0)
)
$r = 22
foreach ($row in $sheet1Rows) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $r++
}

# ---- Sheet4 ("Sheet4", rId4 -> sheet4.xml): append rows 7-9 ----
$ws4 = $wb.Worksheets.Item("Sheet4")
$sheet4Rows = @(
    @("l",    1042.0),
    @("a",    0.0),
    @("Jack", 521.0)
)
$r = 7
foreach ($row in $sheet4Rows) {
    $ws4.Cells.Item($r, 1).Value = $row[0]
    $ws4.Cells.Item($r, 2).Value = $row[1]
    $r++
}

# ---- Sheet5 ("Sheet5", rId5 -> sheet5.xml): append row 7 ----
$ws5 = $wb.Worksheets.Item("Sheet5")
$ws5.Cells.Item(7, 1).Value = "Jack"
$ws5.Cells.Item(7, 2).Value = 16151.0
